$wb = $excel.ActiveWorkbook

# Update column F ("想去人数" / want-to-go count) values on each sheet
# per the refreshed data snapshot (gh-pages output regenerated at 456a3b4).

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 179
$ws.Cells.Item(3, 6).Value = 959
$ws.Cells.Item(4, 6).Value = 1113
$ws.Cells.Item(5, 6).Value = 1565
$ws.Cells.Item(6, 6).Value = 346
$ws.Cells.Item(7, 6).Value = 715
$ws.Cells.Item(8, 6).Value = 12950
$ws.Cells.Item(9, 6).Value = 2242
$ws.Cells.Item(10, 6).Value = 915
$ws.Cells.Item(11, 6).Value = 289
$ws.Cells.Item(12, 6).Value = 53072
$ws.Cells.Item(13, 6).Value = 1280
$ws.Cells.Item(14, 6).Value = 284
$ws.Cells.Item(15, 6).Value = 291
$ws.Cells.Item(16, 6).Value = 843
$ws.Cells.Item(17, 6).Value = 692
$ws.Cells.Item(18, 6).Value = 348
$ws.Cells.Item(19, 6).Value = 2954
$ws.Cells.Item(20, 6).Value = 831
$ws.Cells.Item(21, 6).Value = 4960
$ws.Cells.Item(22, 6).Value = 4960
$ws.Cells.Item(23, 6).Value = 1221
$ws.Cells.Item(24, 6).Value = 916
$ws.Cells.Item(27, 6).Value = 25
$ws.Cells.Item(28, 6).Value = 14
$ws.Cells.Item(29, 6).Value = 1163
$ws.Cells.Item(31, 6).Value = 7
$ws.Cells.Item(32, 6).Value = 137
$ws.Cells.Item(33, 6).Value = 308
$ws.Cells.Item(34, 6).Value = 34
$ws.Cells.Item(36, 6).Value = 52
$ws.Cells.Item(37, 6).Value = 35
$ws.Cells.Item(38, 6).Value = 4623
$ws.Cells.Item(39, 6).Value = 31
$ws.Cells.Item(40, 6).Value = 4687
$ws.Cells.Item(41, 6).Value = 5629
$ws.Cells.Item(43, 6).Value = 141
$ws.Cells.Item(44, 6).Value = 103
$ws.Cells.Item(48, 6).Value = 61
$ws.Cells.Item(49, 6).Value = 4146

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 83
$ws.Cells.Item(5, 6).Value = 117
$ws.Cells.Item(12, 6).Value = 1086

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 774
$ws.Cells.Item(3, 6).Value = 521
$ws.Cells.Item(4, 6).Value = 125

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 521
$ws.Cells.Item(3, 6).Value = 179
$ws.Cells.Item(4, 6).Value = 959
$ws.Cells.Item(5, 6).Value = 1113
$ws.Cells.Item(6, 6).Value = 346
$ws.Cells.Item(7, 6).Value = 715
$ws.Cells.Item(8, 6).Value = 12950
$ws.Cells.Item(9, 6).Value = 12950
$ws.Cells.Item(10, 6).Value = 2242
$ws.Cells.Item(11, 6).Value = 289
$ws.Cells.Item(12, 6).Value = 1280
$ws.Cells.Item(13, 6).Value = 291
$ws.Cells.Item(14, 6).Value = 843
$ws.Cells.Item(15, 6).Value = 692
$ws.Cells.Item(16, 6).Value = 348
$ws.Cells.Item(17, 6).Value = 2954
$ws.Cells.Item(18, 6).Value = 831
$ws.Cells.Item(19, 6).Value = 83
$ws.Cells.Item(20, 6).Value = 4960
$ws.Cells.Item(21, 6).Value = 4960
$ws.Cells.Item(22, 6).Value = 1221
$ws.Cells.Item(24, 6).Value = 117
$ws.Cells.Item(25, 6).Value = 916
$ws.Cells.Item(27, 6).Value = 14
$ws.Cells.Item(28, 6).Value = 1163
$ws.Cells.Item(31, 6).Value = 137
$ws.Cells.Item(33, 6).Value = 308
$ws.Cells.Item(34, 6).Value = 34
$ws.Cells.Item(35, 6).Value = 35
$ws.Cells.Item(36, 6).Value = 4623
$ws.Cells.Item(37, 6).Value = 31
$ws.Cells.Item(38, 6).Value = 4687
$ws.Cells.Item(40, 6).Value = 141
$ws.Cells.Item(41, 6).Value = 103
$ws.Cells.Item(47, 6).Value = 4146
